# This edit re-orders the data rows (rows 2-10) of the "Artfynd" sheet.
# The header row (row 1) is untouched. Each data row's full contents
# (all columns) move as a unit to a new row position per the mapping below,
# which was derived by comparing the before/after column-A "Id" values.
#
# before-row -> after-row
#   2  -> 4
#   3  -> 5
#   4  -> 6
#   5  -> 2
#   6  -> 7
#   7  -> 8
#   8  -> 3
#   9  -> 10
#   10 -> 9
#
# Because this permutation has cycles (e.g. 2->4->6->7->8->3->5->2), we must
# capture every source row's values BEFORE writing any of them, otherwise an
# earlier write would clobber data that a later step still needs to read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 10
$lastCol = 51   # column AY

$mapping = @{
    2  = 4
    3  = 5
    4  = 6
    5  = 2
    6  = 7
    7  = 8
    8  = 3
    9  = 10
    10 = 9
}

# Columns Y/Z/AA/AB hold plain-text date/time strings ("2019-03-12", "00:00")
# that are identical on every data row. Round-tripping them through Value2
# would auto-coerce the text into a real date serial number, silently
# changing the cell's stored type even though the visible content never
# actually changes across this permutation. Skip them - nothing to move.
$skipCols = @(25, 26, 27, 28)

# 1) Snapshot every cell in rows 2-10 (all used columns) into memory, and
#    also remember the row's CURRENT (pre-edit) contents so we can tell
#    whether a destination cell truly needs touching.
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowVals = @{}
    for ($c = 1; $c -le $lastCol; $c++) {
        if ($skipCols -contains $c) { continue }
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write the snapshotted rows back out to their new row positions. Only
#    touch a cell when it actually needs to change value (either the
#    incoming value is non-blank, or the destination currently holds a
#    stale non-blank value that must be cleared because the row moving
#    into that slot has nothing there). This avoids gratuitously rewriting
#    (and thus restructuring) cells that are blank before and after.
foreach ($srcRow in $mapping.Keys) {
    $dstRow = $mapping[$srcRow]
    $rowVals = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        if ($skipCols -contains $c) { continue }
        $newVal = $rowVals[$c]
        $newIsBlank = ($newVal -eq $null) -or ($newVal -eq "")
        $curVal = $ws.Cells.Item($dstRow, $c).Value2
        $curIsBlank = ($curVal -eq $null) -or ($curVal -eq "")
        if ($newIsBlank -and $curIsBlank) { continue }
        $ws.Cells.Item($dstRow, $c).Value2 = $newVal
    }
}
